# The workbook's "Artfynd" sheet lists species-observation records (one per
# row). Rows 2 and 4 turn out to describe the same two observations but with
# their row order swapped (every field - id, taxon info, activity, coords,
# times, biotope text, ... - for row 2 moves to row 4 and vice versa). On
# top of that the easting/northing (Ost/Nord, columns Q/R) for every data
# row (2, 3 and 4) get rounded from a long decimal to a plain integer metre
# value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 4
$lastCol = 51   # column AY

$rowA = 2
$rowB = 4

# Columns whose values look like dates ("2023-08-25") - Excel/COM auto
# parses such literals into date serials when assigned back through
# .Value, even though the source file stores them as plain text in a
# "General" formatted cell. Use the standard Excel "force text" leading
# apostrophe for those specific columns so the round-trip keeps the
# original string representation (and keeps NumberFormat = General).
$textColumns = @(25, 27)   # Y = Startdatum, AA = Slutdatum

# --- Swap every cell of row 2 with the matching cell of row 4 -------------
for ($col = 1; $col -le $lastCol; $col++) {
    $cellA = $ws.Cells.Item($rowA, $col)
    $cellB = $ws.Cells.Item($rowB, $col)

    $valueA = $cellA.Value2
    $valueB = $cellB.Value2

    if ($textColumns -contains $col) {
        if ($valueB -ne $null -and $valueB -ne "") { $valueB = "'" + $valueB }
        if ($valueA -ne $null -and $valueA -ne "") { $valueA = "'" + $valueA }
    }

    $cellA.Value = $valueB
    $cellB.Value = $valueA
}

# --- Round the Ost/Nord (Q/R) coordinates on every data row ---------------
$coordCols = @(17, 18)   # Q = Ost, R = Nord

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    foreach ($col in $coordCols) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = [math]::Round([double]$cell.Value2)
    }
}
